$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values for row 2 (No. Of Monitors per day), keyed by day-of-month.
# Day 22 -> 23 monitors, Day 24 -> 161 monitors, everything else -> 0.
$counts = @{}
for ($day = 1; $day -le 30; $day++) { $counts[$day] = 0 }
$counts[22] = 23
$counts[24] = 161

for ($day = 1; $day -le 30; $day++) {
    $col = $day + 2  # column C (3) corresponds to April 1st

    $dateCell = $ws.Cells.Item(1, $col)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "04/{0:D2}/2020" -f $day
    $dateCell.ClearFormats()

    $ws.Cells.Item(2, $col).Value = $counts[$day]
}
